$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3133.2666
$ws.Range("J17").Value = 3666.3333
$ws.Range("L17").Value = 10998.9999
$ws.Range("N17").Value = -11334.9999

$ws.Range("H64").Value = 1500
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 1500
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 1500
$ws.Range("M64").Value = $null
$ws.Range("N64").Value = -1996

$ws.Range("H67").Value = 1500
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 1500
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 1500
$ws.Range("M67").Value = $null
$ws.Range("N67").Value = -3216

$ws.Range("H132").Value = 847.25
$ws.Range("I132").Value = 847.25
$ws.Range("K132").Value = 2541.75
$ws.Range("M132").Value = -11.75

$ws.Range("H138").Value = 3406
$ws.Range("J138").Value = 3406
$ws.Range("L138").Value = 10218
$ws.Range("N138").Value = -20498

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H50").Value = 22198.75
$ws.Range("I50").Value = 397.5
$ws.Range("J50").Value = 44000
$ws.Range("K50").Value = 397.5
$ws.Range("L50").Value = 44000
$ws.Range("M50").Value = 316.5
$ws.Range("N50").Value = -45428

$ws.Range("H74").Value = 3013.1428
$ws.Range("I74").Value = 3117
$ws.Range("K74").Value = 3117
$ws.Range("M74").Value = -2243

$ws.Range("H77").Value = 3013.1428
$ws.Range("I77").Value = 3117
$ws.Range("K77").Value = 15585
$ws.Range("M77").Value = -11217

$ws.Range("H104").Value = 22916.666
$ws.Range("J104").Value = 22916.666
$ws.Range("L104").Value = 22916.666
$ws.Range("N104").Value = -29904.666

$ws.Range("H122").Value = 8999.5
$ws.Range("J122").Value = 2999
$ws.Range("L122").Value = 8997
$ws.Range("N122").Value = -13897

$ws.Range("H131").Value = 48999
$ws.Range("J131").Value = 48999
$ws.Range("L131").Value = 48999
$ws.Range("N131").Value = -59079

$ws.Range("H132").Value = 2111.9
$ws.Range("I132").Value = 1791
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 5373
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -2843
$ws.Range("N132").Value = -20060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = $null
$ws.Range("N8").Value = $null

$ws.Range("H88").Value = 25312.285
$ws.Range("J88").Value = 25312.285
$ws.Range("L88").Value = 25312.285
$ws.Range("N88").Value = -26124.285

$ws.Range("H91").Value = 25312.285
$ws.Range("J91").Value = 25312.285
$ws.Range("L91").Value = 25312.285
$ws.Range("N91").Value = -28120.285

$ws.Range("H106").Value = 4966.6665
$ws.Range("J106").Value = 4966.6665
$ws.Range("L106").Value = 4966.6665
$ws.Range("N106").Value = -7490.6665

$ws.Range("H134").Value = 1567.0834
$ws.Range("I134").Value = 1618.6364
$ws.Range("K134").Value = 4855.9092
$ws.Range("M134").Value = -2320.9092

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1695
$ws.Range("I31").Value = 1695
$ws.Range("K31").Value = 1695
$ws.Range("M31").Value = -1400

$ws.Range("H34").Value = 1695
$ws.Range("I34").Value = 1695
$ws.Range("K34").Value = 1695
$ws.Range("M34").Value = -1493

$ws.Range("H58").Value = 536.75
$ws.Range("I58").Value = 536.75
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 536.75
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -333.75
$ws.Range("N58").Value = $null

$ws.Range("H62").Value = 3639.4
$ws.Range("I62").Value = 3639.4
$ws.Range("K62").Value = 3639.4
$ws.Range("M62").Value = -3015.4

$ws.Range("H65").Value = 3639.4
$ws.Range("I65").Value = 3639.4
$ws.Range("K65").Value = 18197
$ws.Range("M65").Value = -15077

$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").Value = $null

$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").Value = $null

$ws.Range("H134").Value = 4476.615
$ws.Range("I134").Value = 1472.3636
$ws.Range("J134").Value = 21000
$ws.Range("K134").Value = 4417.0908
$ws.Range("L134").Value = 63000
$ws.Range("M134").Value = -1882.0908
$ws.Range("N134").Value = -68070

$ws.Range("H136").Value = 536.75
$ws.Range("I136").Value = 536.75
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 1610.25
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = 939.75
$ws.Range("N136").Value = $null

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = $null
$ws.Range("N107").Value = $null

$ws.Range("H122").Value = 649.6667
$ws.Range("I122").Value = 674.5
$ws.Range("K122").Value = 6070.5
$ws.Range("M122").Value = -3620.5

$ws.Range("H131").Value = 2015
$ws.Range("I131").Value = 1030
$ws.Range("K131").Value = 3090
$ws.Range("M131").Value = 1950

$ws.Range("H132").Value = 1801.3334
$ws.Range("I132").Value = 1500
$ws.Range("J132").Value = 1952
$ws.Range("K132").Value = 13500
$ws.Range("L132").Value = 17568
$ws.Range("M132").Value = -10970
$ws.Range("N132").Value = -22628

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 2750.75
$ws.Range("I5").Value = 1999.3334
$ws.Range("J5").Value = 5005
$ws.Range("K5").Value = 1999.3334
$ws.Range("L5").Value = 5005
$ws.Range("M5").Value = -1887.3334
$ws.Range("N5").Value = -5229

$ws.Range("H104").Value = 40000
$ws.Range("J104").Value = 40000
$ws.Range("L104").Value = 40000
$ws.Range("N104").Value = -46988

$ws.Range("H132").Value = 2357.2
$ws.Range("I132").Value = 1396.75
$ws.Range("J132").Value = 2997.5
$ws.Range("K132").Value = 4190.25
$ws.Range("L132").Value = 8992.5
$ws.Range("M132").Value = -1660.25
$ws.Range("N132").Value = -14052.5

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 2060.5
$ws.Range("I55").Value = 2412.6
$ws.Range("J55").Value = 300
$ws.Range("K55").Value = 2412.6
$ws.Range("L55").Value = 300
$ws.Range("M55").Value = -2239.6
$ws.Range("N55").Value = -646

$ws.Range("H68").Value = 2000
$ws.Range("I68").Value = 2000
$ws.Range("K68").Value = 2000
$ws.Range("M68").Value = -1251

$ws.Range("H71").Value = 2000
$ws.Range("I71").Value = 2000
$ws.Range("K71").Value = 10000
$ws.Range("M71").Value = -6256

$ws.Range("H136").Value = 1368650.8
$ws.Range("I136").Value = 837500
$ws.Range("K136").Value = 2512500
$ws.Range("M136").Value = -2509950

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 891.5625
$ws.Range("I132").Value = 914.375
$ws.Range("J132").Value = 868.75
$ws.Range("K132").Value = 2743.125
$ws.Range("L132").Value = 2606.25
$ws.Range("M132").Value = -213.125
$ws.Range("N132").Value = -7666.25

$ws.Range("H136").Value = 1452.3462
$ws.Range("I136").Value = 1337.4783
$ws.Range("J136").Value = 2333
$ws.Range("K136").Value = 4012.4349
$ws.Range("L136").Value = 6999
$ws.Range("M136").Value = -1462.4349
$ws.Range("N136").Value = -12099
